# Re-synced cryptos list (prices + 1h volume deltas), plus a Kaspa/Celestia row swap.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# cellRef -> new text value, in the exact order the source rows were re-synced.
$updates = [ordered]@{
    "D2" = "40.360.82"
    "E2" = "  -3.12%  "
    "D3" = "2.365.85"
    "E3" = "  -4.41%  "
    "E4" = "  -0.05%  "
    "D5" = "310.47"
    "E5" = "  -2.79%  "
    "D6" = "86.05"
    "E6" = "  -6.83%  "
    "D7" = "0.532"
    "E7" = "  -3.39%  "
    "E8" = "  +0.01%  "
    "D9" = "0.489"
    "D10" = "0.0831"
    "E10" = "  -3.75%  "
    "D11" = "30.36"
    "E11" = "  -8.16%  "
    "D12" = "0.109"
    "E12" = "  -0.67%  "
    "D13" = "2.737.33"
    "E13" = "  -4.20%  "
    "D14" = "6.46"
    "E14" = "  -6.00%  "
    "D15" = "14.97"
    "E15" = "  -3.53%  "
    "D16" = "2.372.44"
    "E16" = "  -3.97%  "
    "D17" = "0.757"
    "E17" = "  -4.59%  "
    "D18" = "40.357.98"
    "E18" = "  -2.99%  "
    "E19" = "  -3.46%  "
    "D20" = "6.12"
    "E20" = "  -5.07%  "
    "D21" = "68.28"
    "E21" = "  -3.40%  "
    "D22" = "10.82"
    "E22" = "  -3.61%  "
    "D23" = "234.84"
    "E23" = "  -2.03%  "
    "D24" = "2.57"
    "E24" = "  -6.35%  "
    "E25" = "  +0.03%  "
    "D26" = "1.80"
    "E26" = "  -7.26%  "
    "D27" = "23.72"
    "E27" = "  -4.92%  "
    "E28" = "  -4.05%  "
    "E29" = "  -5.15%  "
    "D30" = "34.57"
    "E30" = "  -5.63%  "
    "D31" = "152.77"
    "E31" = "  -2.66%  "
    "E32" = "  +0.00%  "
    "D33" = "5.22"
    "E33" = "  -3.78%  "
    "D34" = "0.0730"
    "E34" = "  -4.41%  "
    "E35" = "  -5.18%  "
    "E36" = "  -1.93%  "
    "E37" = "  -3.31%  "
    "B38" = "Kaspa"
    "C38" = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
    "D38" = "0.0999"
    "E38" = "  -3.88%  "
    "B39" = "Celestia"
    "C39" = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
    "D39" = "15.85"
    "E39" = "  -7.75%  "
    "D40" = "1.72"
    "E40" = "  -6.73%  "
    "D41" = "3.84"
    "E41" = "  -4.20%  "
    "E42" = "  -3.54%  "
    "D43" = "1.968.07"
    "E43" = "  -1.65%  "
    "E44" = "  -5.86%  "
    "D45" = "17.69"
    "E45" = "  -5.29%  "
    "D46" = "9.39"
    "E46" = "  -1.70%  "
    "D47" = "2.68"
    "E47" = "  -9.80%  "
    "D48" = "2.596.52"
    "E48" = "  -4.47%  "
    "D49" = "92.93"
    "E49" = "  -4.78%  "
    "D50" = "71.39"
    "E50" = "  -5.56%  "
    "D51" = "50.40"
    "E51" = "  -3.64%  "
}

# Cells whose new text looks like a plain number (e.g. "310.47") need the column
# forced to Text first, otherwise Excel auto-converts the assignment to a numeric
# value and silently drops meaningful trailing/leading zeros (e.g. "0.0730" -> 0.073).
$forceText = @(
    "D5", "D6", "D7", "D9", "D10", "D11", "D12", "D14", "D15", "D17", "D20", "D21", "D22", "D23", "D24", "D26", "D27", "D30", "D31", "D33", "D34", "D38", "D39", "D40", "D41", "D45", "D46", "D47", "D49", "D50", "D51"
)

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    if ($forceText -contains $ref) {
        $cell.NumberFormat = "@"
        $cell.Value = $updates[$ref]
        $cell.ClearFormats()
    } else {
        $cell.Value = $updates[$ref]
    }
}
